# Apply "adding new progress as of date 04 nov 2025" update:
# For rows 3-35 on the "Training Dashboard" sheet, the PERIOD TO EXPIRE
# (column H) decreases by 1 day and the LAST UPDATE date (column I)
# moves from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 35; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $iCell = $ws.Cells.Item($row, 9)   # column I

    # Keep a copy of H's current formatting (style stays the same for this
    # cell - it is reused below to restore I's formatting after the write).
    $hCell.Copy()
    $current = $hCell.Value2
    $hCell.Value = $current - 1

    # Leading apostrophe keeps this a literal text value ("04-Nov-2025")
    # instead of letting Excel auto-convert the date-looking string into
    # a real date serial number. Writing the apostrophe form flips Excel's
    # "quote prefix" flag on, which would otherwise fork the cell onto a
    # new style record, so re-paste the original (pre-write) formatting
    # from H afterwards to land back on the same style the cell started on.
    $iCell.Value = "'04-Nov-2025"
    $iCell.PasteSpecial(-4122) # xlPasteFormats
}

$excel.CutCopyMode = 0
